# 141: new scrape run 31/12/2025 14:19:08 — appends latest rows to the
# LP1912 sheet and the 6203-6173 sheet, and refreshes the "last updated"
# stamp on all three sheets (LP1912, LP1912-215, 6203-6173).

$wb = $excel.ActiveWorkbook

$newStamp = "Última actualización: 31/12/2025 14:19:08"

# ---------------------------------------------------------------------
# Sheet 1: LP1912  (columns: A Fecha-marker / B Hora_Scrap / C Hora_Llegada /
#                    D Línea / E Minutos / F Parada / G Fecha)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Cells.Item(2, 1).Value = $newStamp
$ws1.Cells.Item(3, 1).Value = "Total filas: 1001"

# Column A on the new rows is left blank (scraper only stamps it when the
# date rolls over). Stamp the blank-cell formatting onto A991:A1002 so the
# cells exist (matching the source feed) without giving them a value.
$ws1.Range("G2").Copy()
$ws1.Range("A991:A1002").PasteSpecial(-4122)

$rows1 = @(
    @("14:18:57", "14:25", "16_SANTA ANA", 7, "LP1912", "31/12/2025"),
    @("14:18:57", "14:33", "23_HERNANDEZ", 15, "LP1912", "31/12/2025"),
    @("14:18:57", "14:37", "16_P MOR-SANTA ANA", 19, "LP1912", "31/12/2025"),
    @("14:18:57", "14:40", "17X38_ROMERO", 22, "LP1912", "31/12/2025"),
    @("14:18:57", "14:49", "16_SANTA ANA", 31, "LP1912", "31/12/2025"),
    @("14:18:57", "15:13", "15_ABASTO", 55, "LP1912", "31/12/2025"),
    @("14:18:57", "15:14", "10_OLMOS", 56, "LP1912", "31/12/2025"),
    @("14:18:57", "15:19", "14_ABASTO", 61, "LP1912", "31/12/2025"),
    @("14:18:57", "15:24", "11_ETCHEVERRY", 66, "LP1912", "31/12/2025"),
    @("14:18:57", "15:33", "16_SANTA ANA", 75, "LP1912", "31/12/2025"),
    @("14:18:57", "15:36", "23_HERNANDEZ", 78, "LP1912", "31/12/2025"),
    @("14:18:57", "15:44", "14_ABASTO", 86, "LP1912", "31/12/2025")
)

$r = 991
foreach ($row in $rows1) {
    $ws1.Cells.Item($r, 2).Value = $row[0]
    $ws1.Cells.Item($r, 3).Value = $row[1]
    $ws1.Cells.Item($r, 4).Value = $row[2]
    $ws1.Cells.Item($r, 5).Value = $row[3]
    $ws1.Cells.Item($r, 6).Value = $row[4]
    $ws1.Cells.Item($r, 7).Value = $row[5]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Sheet 2: LP1912-215 — only the "last updated" stamp changes.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2, 1).Value = $newStamp

# ---------------------------------------------------------------------
# Sheet 3: 6203-6173  (columns: A Fecha-marker / B Fecha / C Hora_Scrap /
#                       D Hora_Llegada / E Línea / F Minutos / G Parada)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Cells.Item(2, 1).Value = $newStamp
$ws3.Cells.Item(3, 1).Value = "Total filas: 123"

# Same blank-but-present column A treatment for the two appended rows.
$ws3.Range("A4").Copy()
$ws3.Range("A123:A124").PasteSpecial(-4122)

$rows3 = @(
    @("31/12/2025", "14:19:03", "14:35", "215C_LA PLATA", 16, "L6203"),
    @("31/12/2025", "14:19:08", "15:01", "215A_LA PLATA", 42, "L6173")
)

$r = 123
foreach ($row in $rows3) {
    $ws3.Cells.Item($r, 2).Value = $row[0]
    $ws3.Cells.Item($r, 3).Value = $row[1]
    $ws3.Cells.Item($r, 4).Value = $row[2]
    $ws3.Cells.Item($r, 5).Value = $row[3]
    $ws3.Cells.Item($r, 6).Value = $row[4]
    $ws3.Cells.Item($r, 7).Value = $row[5]
    $r = $r + 1
}
